$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for specific rows to reflect repulled data / recalculated means
$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -9
$ws.Range("F5").Value = -4
$ws.Range("F7").Value = -2
$ws.Range("F10").Value = -5
